$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 2).Value = 0.1351351351351351
$ws.Cells.Item(2, 3).Value = 0.6418918918918919
$ws.Cells.Item(2, 10).Value = 0.01576576576576576
$ws.Cells.Item(2, 16).Value = 0.0990990990990991
$ws.Cells.Item(2, 19).Value = 0.1081081081081081

# Row 3
$ws.Cells.Item(3, 2).Value = 0.006688963210702341
$ws.Cells.Item(3, 3).Value = 0.03678929765886288
$ws.Cells.Item(3, 10).Value = 0.01337792642140468
$ws.Cells.Item(3, 16).Value = 0.7892976588628763
$ws.Cells.Item(3, 19).Value = 0.1538461538461539

# Row 4
$ws.Cells.Item(4, 10).Value = 0.06521739130434782
$ws.Cells.Item(4, 16).Value = 0.6086956521739131
$ws.Cells.Item(4, 19).Value = 0.3260869565217391

# Row 5
$ws.Cells.Item(5, 16).Value = 0.8
$ws.Cells.Item(5, 19).Value = 0.2

# Row 6
$ws.Cells.Item(6, 2).Value = 0.05034324942791762
$ws.Cells.Item(6, 4).Value = 0.011441647597254
$ws.Cells.Item(6, 6).Value = 0.06636155606407322
$ws.Cells.Item(6, 10).Value = 0.2334096109839817
$ws.Cells.Item(6, 15).Value = 0.0137299771167048
$ws.Cells.Item(6, 17).Value = 0.1876430205949657
$ws.Cells.Item(6, 18).Value = 0.06864988558352403
$ws.Cells.Item(6, 19).Value = 0.3684210526315789

# Row 7
$ws.Cells.Item(7, 2).Value = 0.1016042780748663
$ws.Cells.Item(7, 4).Value = 0.0374331550802139
$ws.Cells.Item(7, 6).Value = 0.04545454545454546
$ws.Cells.Item(7, 10).Value = 0.1336898395721925
$ws.Cells.Item(7, 15).Value = 0.02406417112299465
$ws.Cells.Item(7, 17).Value = 0.1550802139037433
$ws.Cells.Item(7, 18).Value = 0.09090909090909091
$ws.Cells.Item(7, 19).Value = 0.4117647058823529

# Row 8
$ws.Cells.Item(8, 2).Value = 0.07516339869281045
$ws.Cells.Item(8, 4).Value = 0.02287581699346405
$ws.Cells.Item(8, 5).Value = 0.001089324618736384
$ws.Cells.Item(8, 6).Value = 0.08169934640522876
$ws.Cells.Item(8, 10).Value = 0.09586056644880174
$ws.Cells.Item(8, 15).Value = 0.01416122004357298
$ws.Cells.Item(8, 17).Value = 0.1840958605664488
$ws.Cells.Item(8, 18).Value = 0.1045751633986928
$ws.Cells.Item(8, 19).Value = 0.420479302832244

# Row 9
$ws.Cells.Item(9, 2).Value = 0.07861635220125786
$ws.Cells.Item(9, 4).Value = 0.01572327044025157
$ws.Cells.Item(9, 5).Value = 0.003144654088050315
$ws.Cells.Item(9, 6).Value = 0.09119496855345911
$ws.Cells.Item(9, 10).Value = 0.07232704402515723
$ws.Cells.Item(9, 15).Value = 0.02515723270440252
$ws.Cells.Item(9, 17).Value = 0.1886792452830189
$ws.Cells.Item(9, 18).Value = 0.1226415094339623
$ws.Cells.Item(9, 19).Value = 0.4025157232704403

# Row 10
$ws.Cells.Item(10, 2).Value = 0.1014686248331108
$ws.Cells.Item(10, 4).Value = 0.02225189141076991
$ws.Cells.Item(10, 5).Value = 0.001335113484646195
$ws.Cells.Item(10, 6).Value = 0.07610146862483311
$ws.Cells.Item(10, 10).Value = 0.09479305740987984
$ws.Cells.Item(10, 15).Value = 0.01424121050289275
$ws.Cells.Item(10, 17).Value = 0.225634178905207
$ws.Cells.Item(10, 18).Value = 0.08900756564307966
$ws.Cells.Item(10, 19).Value = 0.3751668891855808

# Row 11
$ws.Cells.Item(11, 7).Value = 0.1686956521739131
$ws.Cells.Item(11, 10).Value = 0.07478260869565218
$ws.Cells.Item(11, 11).Value = 0.1965217391304348
$ws.Cells.Item(11, 12).Value = 0.5478260869565217
$ws.Cells.Item(11, 19).Value = 0.01217391304347826

# Row 12
$ws.Cells.Item(12, 7).Value = 0.7376093294460642
$ws.Cells.Item(12, 10).Value = 0.1836734693877551
$ws.Cells.Item(12, 11).Value = 0.008746355685131196
$ws.Cells.Item(12, 12).Value = 0.03206997084548105
$ws.Cells.Item(12, 19).Value = 0.03790087463556852

# Row 13
$ws.Cells.Item(13, 6).Value = 0.02040816326530612
$ws.Cells.Item(13, 7).Value = 0.5
$ws.Cells.Item(13, 10).Value = 0.3877551020408163
$ws.Cells.Item(13, 19).Value = 0.09183673469387756

# Row 14
$ws.Cells.Item(14, 7).Value = 0.3333333333333333
$ws.Cells.Item(14, 10).Value = 0.3333333333333333
$ws.Cells.Item(14, 19).Value = 0.3333333333333333

# Row 15
$ws.Cells.Item(15, 6).Value = 0.02117647058823529
$ws.Cells.Item(15, 8).Value = 0.1835294117647059
$ws.Cells.Item(15, 9).Value = 0.06588235294117648
$ws.Cells.Item(15, 11).Value = 0.0611764705882353
$ws.Cells.Item(15, 13).Value = 0.01411764705882353
$ws.Cells.Item(15, 15).Value = 0.05411764705882353
$ws.Cells.Item(15, 19).Value = 0.2235294117647059

# Row 16
$ws.Cells.Item(16, 6).Value = 0.00911854103343465
$ws.Cells.Item(16, 8).Value = 0.1945288753799392
$ws.Cells.Item(16, 9).Value = 0.05167173252279635
$ws.Cells.Item(16, 10).Value = 0.4164133738601823
$ws.Cells.Item(16, 11).Value = 0.1155015197568389
$ws.Cells.Item(16, 13).Value = 0.0243161094224924
$ws.Cells.Item(16, 14).Value = 0.00303951367781155
$ws.Cells.Item(16, 15).Value = 0.06382978723404255
$ws.Cells.Item(16, 19).Value = 0.121580547112462

# Row 17
$ws.Cells.Item(17, 6).Value = 0.01487414187643021
$ws.Cells.Item(17, 8).Value = 0.2162471395881007
$ws.Cells.Item(17, 9).Value = 0.07894736842105263
$ws.Cells.Item(17, 10).Value = 0.4221967963386727
$ws.Cells.Item(17, 11).Value = 0.08009153318077804
$ws.Cells.Item(17, 13).Value = 0.01830663615560641
$ws.Cells.Item(17, 14).Value = 0.0011441647597254
$ws.Cells.Item(17, 15).Value = 0.06636155606407322
$ws.Cells.Item(17, 19).Value = 0.1018306636155606

# Row 18
$ws.Cells.Item(18, 6).Value = 0.007556675062972292
$ws.Cells.Item(18, 8).Value = 0.2241813602015113
$ws.Cells.Item(18, 9).Value = 0.07556675062972293
$ws.Cells.Item(18, 10).Value = 0.4030226700251889
$ws.Cells.Item(18, 11).Value = 0.1007556675062972
$ws.Cells.Item(18, 13).Value = 0.03526448362720403
$ws.Cells.Item(18, 15).Value = 0.07052896725440806
$ws.Cells.Item(18, 19).Value = 0.08312342569269521

# Row 19
$ws.Cells.Item(19, 6).Value = 0.02024647887323944
$ws.Cells.Item(19, 8).Value = 0.2143485915492958
$ws.Cells.Item(19, 9).Value = 0.07614436619718309
$ws.Cells.Item(19, 10).Value = 0.3631161971830986
$ws.Cells.Item(19, 11).Value = 0.1232394366197183
$ws.Cells.Item(19, 13).Value = 0.02332746478873239
$ws.Cells.Item(19, 14).Value = 0.0008802816901408451
$ws.Cells.Item(19, 15).Value = 0.07614436619718309
$ws.Cells.Item(19, 19).Value = 0.1025528169014085
